$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).ClearFormats()
}

Set-PriceText "D2" "63.120.53"
$ws.Range("E2").Value = "  +0.09%  "

Set-PriceText "D3" "2.565.62"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-PriceText "D5" "585.87"
$ws.Range("E5").Value = "  +0.15%  "

Set-PriceText "D6" "143.80"
$ws.Range("E6").Value = "  -2.78%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -1.27%  "

$ws.Range("E9").Value = "  -2.49%  "

Set-PriceText "D10" "5.62"
$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("E12").Value = "  -1.85%  "

Set-PriceText "D13" "27.03"
$ws.Range("E13").Value = "  -1.23%  "

Set-PriceText "D14" "3.026.21"
$ws.Range("E14").Value = "  +0.14%  "

Set-PriceText "D15" "63.017.21"
$ws.Range("E15").Value = "  +0.01%  "

$ws.Range("E16").Value = "  -1.07%  "

Set-PriceText "D17" "2.566.77"
$ws.Range("E17").Value = "  -0.82%  "

$ws.Range("E18").Value = "  -2.99%  "

Set-PriceText "D19" "339.61"
$ws.Range("E19").Value = "  -1.28%  "

$ws.Range("E20").Value = "  -2.25%  "

$ws.Range("E21").Value = "  -3.82%  "

$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("E23").Value = "  +3.70%  "

Set-PriceText "D24" "67.62"
$ws.Range("E24").Value = "  +1.66%  "

Set-PriceText "D25" "1.55"
$ws.Range("E25").Value = "  +3.96%  "

$ws.Range("E26").Value = "  -1.24%  "

Set-PriceText "D27" "0.163"
$ws.Range("E27").Value = "  -3.86%  "

$ws.Range("E28").Value = "  -0.30%  "

Set-PriceText "D29" "7.91"
$ws.Range("E29").Value = "  -2.53%  "

Set-PriceText "D30" "8.16"
$ws.Range("E30").Value = "  -2.89%  "

$ws.Range("E31").Value = "  -1.68%  "

Set-PriceText "D32" "467.37"
$ws.Range("E32").Value = "  +1.22%  "

$sub3 = [char]0x2083
Set-PriceText "D33" "0.0${sub3}0794"
$ws.Range("E33").Value = "  -3.23%  "

$ws.Range("E34").Value = "  +2.69%  "

Set-PriceText "D35" "175.99"
$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("E37").Value = "  -2.67%  "

Set-PriceText "D38" "18.78"
$ws.Range("E38").Value = "  -1.88%  "

Set-PriceText "D39" "4.52"
$ws.Range("E39").Value = "  -0.23%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  -3.42%  "

$ws.Range("E42").Value = "  +1.62%  "

Set-PriceText "D43" "156.90"
$ws.Range("E43").Value = "  +4.12%  "

Set-PriceText "D44" "3.67"
$ws.Range("E44").Value = "  -3.98%  "

Set-PriceText "D45" "21.10"
$ws.Range("E45").Value = "  +1.22%  "

$ws.Range("E46").Value = "  +2.63%  "

$ws.Range("E47").Value = "  -1.72%  "

$ws.Range("E48").Value = "  -1.33%  "

$ws.Range("E49").Value = "  -1.54%  "

Set-PriceText "D50" "17.99"
$ws.Range("E50").Value = "  -2.18%  "

Set-PriceText "D51" "11.39"
$ws.Range("E51").Value = "  -0.06%  "
